$d = $word.ActiveDocument

# 1. Remove the four red "Research Questions" bullet paragraphs that were
#    dropped from the list (numId 6), right before "Background Research Objectives".
$targets = @(
  "How to test each feature after development of the project?",
  "How can you setup deployment environment of the project?",
  "What is the existing payment API",
  "How can a system architecture be implemented to help the development website application?"
)

$paras = @()
foreach ($p in $d.Paragraphs) {
  $t = $p.Range.Text
  foreach ($target in $targets) {
    if ($t -like "*$target*") {
      $paras += $p
    }
  }
}

for ($i = $paras.Count - 1; $i -ge 0; $i--) {
  $paras[$i].Range.Delete()
}

# 2. Remove the stray <w:lastRenderedPageBreak/> marker that precedes
#    "To meet the deadline of the clients." by re-writing that run's text,
#    which regenerates the run without the stale rendering marker.
$d.Content.Find.Execute("To meet the deadline of the clients.", $true, $false, $false, $false, $false, $true, 1, $false, "To meet the deadline of the clients.", 2) | Out-Null

$d.Save()
